{"js": "// Fix wording: \"unbenutzt\" -> \"ungenutzt\" in the Infoschreiben Maxtarif letter.\nconst searchResults = context.document.body.search(\"unbenutzt\", { matchCase: true, matchWholeWord: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < searchResults.items.length; i++) {\n  searchResults.items[i].insertText(\"ungenutzt\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Fix wording: \"unbenutzt\" -> \"ungenutzt\" in the Infoschreiben Maxtarif letter.\n$d = $word.ActiveDocument\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"unbenutzt\"\n$find.Replacement.Text = \"ungenutzt\"\n$find.Execute(\"unbenutzt\", $true, $true, $false, $false, $false, $true, 1, $false, \"ungenutzt\", 2)\n"}
